$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.383.98'
$ws.Range("E2").Value = '''  +0.12%  '

$ws.Range("D3").Value = '''1.866.88'
$ws.Range("E3").Value = '''  -0.54%  '

$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '''  +0.03%  '

$ws.Range("D5").Value = '''243.45'
$ws.Range("E5").Value = '''  +0.24%  '

$ws.Range("D6").Value = '''0.7001'
$ws.Range("E6").Value = '''  -2.75%  '

$ws.Range("E7").Value = '''  -0.06%  '

$ws.Range("D8").Value = '''0.07871'
$ws.Range("E8").Value = '''  -1.77%  '

$ws.Range("D9").Value = '''0.3118'
$ws.Range("E9").Value = '''  -0.41%  '

$ws.Range("D10").Value = '''24.32'
$ws.Range("E10").Value = '''  -2.07%  '

$ws.Range("D11").Value = '''0.07764'
$ws.Range("E11").Value = '''  -4.72%  '

$ws.Range("D12").Value = '''1.879.29'
$ws.Range("E12").Value = '''  -0.19%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''5.143'
$ws.Range("E13").Value = '''  -1.54%  '

$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").Value = '''92.55'
$ws.Range("E14").Value = '''  -2.01%  '

$ws.Range("D15").Value = '''0.6968'
$ws.Range("E15").Value = '''  -1.69%  '

$ws.Range("D16").Value = '''6.497'
$ws.Range("E16").Value = '''  +1.66%  '

$ws.Range("E17").Value = '''  +1.06%  '

$ws.Range("D18").Value = '''29.406.72'
$ws.Range("E18").Value = '''  +0.19%  '

$ws.Range("D19").Value = '''248.59'
$ws.Range("E19").Value = '''  +1.91%  '

$ws.Range("D20").Value = '''2.120.91'
$ws.Range("E20").Value = '''  -0.25%  '

$ws.Range("D21").Value = '''13.00'
$ws.Range("E21").Value = '''  -1.84%  '

$ws.Range("D22").Value = '''0.9997'
$ws.Range("E22").Value = '''  -0.24%  '

$ws.Range("D23").Value = '''7.579'
$ws.Range("E23").Value = '''  -1.97%  '

$ws.Range("D24").Value = '''1.004'
$ws.Range("E24").Value = '''  +0.09%  '

$ws.Range("D25").Value = '''0.1533'
$ws.Range("E25").Value = '''  -4.52%  '

$ws.Range("D26").Value = '''8.960'
$ws.Range("E26").Value = '''  -0.75%  '

$ws.Range("D27").Value = '''160.79'
$ws.Range("E27").Value = '''  -1.01%  '

$ws.Range("D28").Value = '''18.64'
$ws.Range("E28").Value = '''  +0.87%  '

$ws.Range("D29").Value = '''1.588'
$ws.Range("E29").Value = '''  +5.51%  '

$ws.Range("D30").Value = '''4.278'
$ws.Range("E30").Value = '''  -2.76%  '

$ws.Range("D31").Value = '''4.233'
$ws.Range("E31").Value = '''  -0.85%  '

$ws.Range("D32").Value = '''1.203'
$ws.Range("E32").Value = '''  -1.47%  '

$ws.Range("D33").Value = '''0.05238'
$ws.Range("E33").Value = '''  -1.96%  '

$ws.Range("D34").Value = '''1.881'
$ws.Range("E34").Value = '''  -2.60%  '

$ws.Range("D35").Value = '''0.7544'
$ws.Range("E35").Value = '''  -0.94%  '

$ws.Range("D36").Value = '''1.176'
$ws.Range("E36").Value = '''  +0.18%  '

$ws.Range("D37").Value = '''2.705'
$ws.Range("E37").Value = '''  +0.20%  '

$ws.Range("D38").Value = '''0.01860'
$ws.Range("E38").Value = '''  -0.35%  '

$ws.Range("D39").Value = '''1.268.34'
$ws.Range("E39").Value = '''  +0.19%  '

$ws.Range("D40").Value = '''2.745'
$ws.Range("E40").Value = '''  -0.50%  '

$ws.Range("D41").Value = '''0.8973'
$ws.Range("E41").Value = '''  -0.78%  '

$ws.Range("D42").Value = '''109.42'
$ws.Range("E42").Value = '''  -3.05%  '

$ws.Range("D43").Value = '''5.951'
$ws.Range("E43").Value = '''  -7.41%  '

$ws.Range("D44").Value = '''70.21'
$ws.Range("E44").Value = '''  -4.94%  '

$ws.Range("D45").Value = '''1.000'
$ws.Range("E45").Value = '''  -0.18%  '

$ws.Range("E46").Value = '''  -2.55%  '

$ws.Range("D47").Value = '''2.026.33'
$ws.Range("E47").Value = '''  +0.21%  '

$ws.Range("D48").Value = '''9.584'

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '''0.5180'
$ws.Range("E49").Value = '''  -0.32%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '''1.788'
$ws.Range("E50").Value = '''  -0.32%  '

$ws.Range("D51").Value = '''0.4271'
$ws.Range("E51").Value = '''  -1.30%  '
